$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 79) mirrors the existing rows: every column holds its
# value as literal text (dates like "2025-05-19", numbers-as-text like
# "37.5" or "5,298"), matching how the sheet already stores rows 2-78.
$row = 79
$values = @(
    "2025-05-19",
    "37.5",
    "37",
    "0.94",
    "0.258",
    "0.09",
    "5,298",
    "7,931",
    "7,981",
    "7.2226"
)

$range = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 10))
# Force text storage so Excel doesn't auto-coerce these into dates/numbers.
$range.NumberFormat = "@"

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($row, $i + 1).Value = $values[$i]
}

# Drop the explicit "@" format again so the new row ends up styleless,
# just like the rest of the sheet's data rows.
$range.ClearFormats()
